# Adds System test cases REQ57-61 to the "Test Cases & Results" sheet.
# Cell writes are ordered to reproduce the original author's shared-string
# table append order (new unique strings are appended to xl/sharedStrings.xml
# in first-write order).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases & Results")

# --- Row 60 : REQ-57 ------------------------------------------------------
$ws.Range("F60").Value = "Test that There should also be a general information page showing the number of books that’s on loan, reserved, available books,"
$ws.Range("G60").Value = "Website is accessed"
$ws.Range("H60").Value = "On a browser while connecting to the same network as the RPI, type in the RPI's IP address:5000 "
$ws.Range("I60").Value = "general information page showing the number of books that’s on loan, reserved, available books is shown"
$ws.Range("J60").Value = "general information page showing the number of books that’s on loan, reserved, available books is shown"
$ws.Range("E60").Value = "Mid Impact"
$ws.Rows.Item(60).RowHeight = 86.4

# --- Row 61 : REQ-58 ------------------------------------------------------
$ws.Range("F61").Value = "Test that Logs should be added, every time an “event” happens, like when a book is returned, loaned or a problem is detected on site, a log is created and should be displayed in a page on the website"
$ws.Range("G61").Value = "Open logs is clicked"
$ws.Range("H61").Value = "On a browser while connecting to the same network as the RPI, type in the RPI's IP address:5000 "
$ws.Range("I61").Value = "The time, type of log and log msg of every event that has occurred after the RPI is lauched is shown"
$ws.Range("J61").Value = "The time, type of log and log msg of every event that has occurred after the RPI is lauched is shown"
$ws.Range("E61").Value = "High Impact"
$ws.Rows.Item(61).RowHeight = 100.8

# --- Row 62 : REQ-59 ------------------------------------------------------
$ws.Range("F62").Value = "Test that There should be a page for staff to add books into circulation"
$ws.Range("G62").Value = "Website is accessed"
$ws.Range("H62").Value = "On a browser while connecting to the same network as the RPI, type in the RPI's IP address:5000 then click on the new book button"
$ws.Range("I62").Value = "A screen comes up prompting the staff to input book details like id and title and location which is then used to create book on firebase"
$ws.Range("J62").Value = "A screen comes up prompting the staff to input book details like id and title and location which is then used to create book on firebase"
$ws.Range("E62").Value = "Mid Impact"
$ws.Rows.Item(62).RowHeight = 100.8

# --- Row 63 : REQ-60 (H63 text is deferred - see below) -------------------
$ws.Range("F63").Value = "Test that There should be a page for staff to remove books from circulation"
$ws.Range("G63").Value = "Bookedit is accessed"
$ws.Range("I63").Value = "if updated when red, book will be removed from firebasee"
$ws.Range("J63").Value = "if updated when red, book will be removed from firebasee"
$ws.Range("E63").Value = "Mid Impact"
$ws.Rows.Item(63).RowHeight = 72

# --- Row 64 : REQ-61 -------------------------------------------------------
$ws.Range("F64").Value = "Test that There should be a page for staff to remove user accounts"
$ws.Range("G64").Value = "User edit is accessed"
$ws.Range("H64").Value = "There is a button labeled delete when updating user details, pressing it will cause the screen to toggle red and white, signaling delete or not"

# H63 is written after H64 (matches the original shared-string append order).
$ws.Range("H63").Value = "There is a button labeled delete when updating book details, pressing it will cause the screen to toggle red and white, signaling delete or not"

$ws.Range("I64").Value = "if updated when red, user will be removed from firebasee"
$ws.Range("J64").Value = "if updated when red, user will be removed from firebasee"
$ws.Range("E64").Value = "Mid Impact"
$ws.Rows.Item(64).RowHeight = 72

# --- View state: scroll position + active selection ----------------------
$ws.Range("H63").Select()
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 1
